$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update K column (최종점수 / final score) values
$ws.Range("K2").Value = 60.8
$ws.Range("K3").Value = 58.4
$ws.Range("K4").Value = 51.6
$ws.Range("K5").Value = 48.8

# Update N column (MACRO_SCORE) values
$ws.Range("N2").Value = 54.82400714602223
$ws.Range("N3").Value = 54.82400714602223
$ws.Range("N4").Value = 54.82400714602223
$ws.Range("N5").Value = 54.82400714602223
